# Add "Out-of-Hospital Deaths" row (row 7) to the "ethnicities" and "prop"
# sheets, mirroring the source workbook's upload diff.

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("ethnicities")
$ws3 = $wb.Worksheets.Item("prop")

# ---------------------------------------------------------------------------
# ethnicities sheet: new row 7 = row 5 (Hospitalizations) minus row 6 (ICU)
# ---------------------------------------------------------------------------
$ws2.Range("A7").Value   = "Out-of-Hospital Deaths"
$ws2.Range("B7").Formula = "=B5-B6"
$ws2.Range("C7").Formula = "=C5-C6"
$ws2.Range("D7").Formula = "=D5-D6"
$ws2.Range("E7").Formula = "=E5-E6"
$ws2.Range("F7").Formula = "=F5-F6"
$ws2.Range("G7").Formula = "=G5-G6"

# ---------------------------------------------------------------------------
# prop sheet: new row 7 = percentages matching ethnicities row 7
# ---------------------------------------------------------------------------
$ws3.Range("A7").Value = "Out-of-Hospital Deaths"
$ws3.Range("B7").Value = 79.285714285714278
$ws3.Range("C7").Value = 11.428571428571429
$ws3.Range("D7").Value = 4.2857142857142856
$ws3.Range("E7").Value = 2.1428571428571428
$ws3.Range("F7").Value = 0.7142857142857143
$ws3.Range("G7").Value = 2.1428571428571428

# ---------------------------------------------------------------------------
# Widen column A on both sheets to fit the new, longer label.
# (ColumnWidth is quantized to whole pixels, same as real Excel, so we bias
# the input so the stored width lands as close as possible to the target.)
# ---------------------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 19.276041666666668
$ws3.Columns.Item(1).ColumnWidth = 19.721354166666668

# ---------------------------------------------------------------------------
# Selections: extend the ethnicities selection to include the new row, then
# move the prop sheet's active cell to C8 (leaving prop as the active tab).
# ---------------------------------------------------------------------------
[void]$ws2.Range("A1:G7").Select()
[void]$ws3.Range("C8").Select()
